$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) relabeling ---
$ws.Range("O1").Value = "car"
$ws.Range("P1").Value = "Spcl Allowance"
$ws.Range("Q1").Value = "Arrears"
$ws.Range("R1").Value = "Gross Pay"
$ws.Range("S1").Value = "PF"
$ws.Range("T1").Value = "ESIC"
$ws.Range("U1").Value = "PT"
$ws.Range("V1").Value = "TDS"
$ws.Range("W1").Value = "Deductible Arrears"
$ws.Range("X1").Value = "house"

# --- Employee data row (row 2) ---
$ws.Range("C2").Value = "Vidya Sagar  Pogiri"

# D2 holds a literal text date string (not a real date) in the source
# workbook, so force text entry and restore the default "Normal" style
# to avoid Excel auto-converting it to a date serial number.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2014-03-03"
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = "Regular"
$ws.Range("F2").Value = "Senior HR Executive"

$ws.Range("H2").Value = 12000
$ws.Range("I2").Value = 1000
$ws.Range("J2").Value = 5434.5
$ws.Range("M2").Value = 400
$ws.Range("N2").Value = 100
$ws.Range("P2").Value = 434.5
$ws.Range("Q2").Value = 4500
$ws.Range("R2").Value = 5434.5
$ws.Range("S2").Value = 48
$ws.Range("T2").Value = 95.1
$ws.Range("U2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 143.1
$ws.Range("Z2").Value = 5291.4
